$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on every cell we touch so Excel COM does not
# reinterpret numeric-looking strings (e.g. "1.00", "0.533") as numbers,
# and so multi-dot price strings (e.g. "27.804.96") stay verbatim text.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.804.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.649.61'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.41%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.50'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.533'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.95%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '23.16'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -1.66%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.65%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0615'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0890'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.883.85'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.644.06'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.51%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.55%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.32%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '64.49'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -1.57%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '27.776.38'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +1.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '233.68'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.72%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.69'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.50%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0725'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.22%  '
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.31'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.11'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +7.58%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '150.56'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.44%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.34%  '
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.01%  '
$ws.Range('B28').NumberFormat = '@'
$ws.Range('B28').Value = 'BinanceUSD'
$ws.Range('C28').NumberFormat = '@'
$ws.Range('C28').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.17%  '
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.67'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.03%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.21%  '
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.05%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.32%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +1.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.444.30'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.38%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.35'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.82%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.66%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.885'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -2.24%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.874'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +10.29%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.16%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.60'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.91%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.70'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +2.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.47'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.792.49'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.73'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.62%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '86.41'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.72%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +2.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0997'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.32%  '
